$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the name of the matlab package: "BNT" -> "Global MIT" (row 8, column A)
$ws.Range("A8").Value = "Global MIT"

# Move the active selection to B12, matching the author's final cursor position
[void]$ws.Range("B12").Select()

# Column A widens to best-fit the new, longer "Global MIT" text
$ws.Columns.Item(1).ColumnWidth = 9
